$ws = $excel.ActiveWorkbook.ActiveSheet

# Column indices for the data block (B,C,D,E,F,I,J,K,L,M); G stays 1, H is blank
$cols = @(2,3,4,5,6,9,10,11,12,13)

# New values for rows 2-25 (24 data rows), one row per array
$data = @(
    @(1.02,1.045635148374532,1.048989521873808,1.043272374191927,1.057266750859573,1.038945101223043,1.050693830037792,1.051747829631311,1.046046707304128,1.060002189996295),
    @(1.02,1.047219679848313,1.050208786366259,1.044642879433515,1.058716748892199,1.039344011464254,1.051922970927945,1.052778198188464,1.047226751078704,1.061264374463785),
    @(1.02,1.048243119561687,1.050995869573817,1.045528195622351,1.059653462628368,1.03959960859058,1.052716100873802,1.053442493544439,1.047988300916332,1.062079045932167),
    @(1.02,1.048672938301473,1.051326319781069,1.045900032631819,1.060046897896429,1.039706459944936,1.053049011617624,1.053721189366414,1.048307980266619,1.06242105085341),
    @(1.02,1.048745081507456,1.051381778230287,1.045962445395373,1.060112936565062,1.039724365552408,1.053104878445934,1.053767950134936,1.048361628107235,1.062478446863438),
    @(1.02,1.048248864519481,1.051000286785134,1.045533165492534,1.059658721131713,1.039601038704505,1.052720551282191,1.05344621973964,1.047992574350271,1.062083617707869),
    @(1.02,1.046171038684979,1.049401967457019,1.043735856296062,1.057757106253325,1.039080439140627,1.051109685113578,1.052096552556711,1.046445930423231,1.060429179023084),
    @(1.02,1.042494985065503,1.046570976444088,1.040557007290881,1.054394144293653,1.038143620304707,1.048253913603922,1.049699451605135,1.04370479083482,1.057497862877583),
    @(1.02,1.040033825574568,1.044673467546508,1.038429402148724,1.052143572707677,1.037505821568433,1.046338036825167,1.048088360404428,1.041866343152298,1.0555324667126),
    @(1.02,1.03896550157517,1.04384932035441,1.037506035510072,1.051166903744556,1.037226465435273,1.045505486591483,1.047387567017714,1.041067564905601,1.054678676948245),
    @(1.02,1.038568272322504,1.043542809940863,1.037162731520864,1.050803792036691,1.037122218179206,1.045195786394256,1.047126776011952,1.040770446352182,1.054361118131174),
    @(1.02,1.038653497913347,1.043608575034438,1.037236386223024,1.050881695967521,1.037144601434201,1.045262238823524,1.047182738640142,1.040834198284155,1.054429254868614),
    @(1.02,1.038932674834492,1.043823992016807,1.037477664567828,1.051136895672603,1.037217858176168,1.04547989598877,1.047366019891234,1.041043013543421,1.054652436122184),
    @(1.02,1.039104630779597,1.043956666216999,1.037626280846236,1.051294088177392,1.037262930106416,1.045613941330967,1.047478881012616,1.041171616042079,1.054789889110966),
    @(1.02,1.040104670405722,1.044728109860073,1.038490637943634,1.052208344775701,1.037524294117333,1.046393227265235,1.048134802124072,1.041919297451356,1.055589071040923),
    @(1.02,1.040731257661335,1.045211337668915,1.039032257821975,1.05278124981637,1.037687385712982,1.046881253275848,1.048545387500966,1.042387564611167,1.056089631716148),
    @(1.02,1.041096482858632,1.045492954310588,1.039347973572573,1.053115208189176,1.037782207202128,1.047165625310371,1.048784568655291,1.042660435218415,1.056381334549297),
    @(1.02,1.041220972727471,1.045588937502713,1.039455590501553,1.053229044459482,1.037814486896561,1.047262540702198,1.048866071506495,1.042753432905747,1.056480752861326),
    @(1.02,1.040664056985425,1.045159516982541,1.038974168090296,1.052719804052166,1.037669919315727,1.046828922275313,1.048501367316058,1.04233735110639,1.056035953841831),
    @(1.02,1.038850475495544,1.043760567788422,1.037406623140204,1.05106175499837,1.037196299223432,1.045415813981647,1.047312061582469,1.040981534225903,1.054586726559714),
    @(1.02,1.037707849105032,1.042878758857783,1.036419163366945,1.050017340280796,1.036895725334611,1.044524705745159,1.046561487311233,1.04012666297673,1.053673086547201),
    @(1.02,1.038313804192128,1.043346436859715,1.036942815921069,1.050571190701787,1.037055330867217,1.044997351432843,1.046959649810416,1.04058007815864,1.054157659836256),
    @(1.02,1.040694422858709,1.045182933267954,1.039000416963806,1.052747569389285,1.037677812577749,1.046852569283523,1.048521259089531,1.042360041236564,1.056060209383978),
    @(1.02,1.043447131796739,1.047304620475771,1.041380257306681,1.05526503043806,1.03838813362123,1.048994285846524,1.050321427589336,1.04441535208128,1.058257617455097)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $vals[$j]
    }
}
